$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04217225448213213
$ws.Range("D2").Value = 0.6254007317898527
$ws.Range("E2").Value = 0.07971969209722118
$ws.Range("F2").Value = 7.677665792954542
$ws.Range("G2").Value = 0.002619586479037079
$ws.Range("I2").Value = 6.702089430277283
$ws.Range("L2").Value = 0.230991611700361
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("C3").Value = 0.03922197126733096
$ws.Range("D3").Value = 0.5997560311011512
$ws.Range("E3").Value = 0.07954608488608628
$ws.Range("F3").Value = 7.312664915337109
$ws.Range("G3").Value = 0.00263784524586902
$ws.Range("I3").Value = 6.40758621250211
$ws.Range("L3").Value = 0.2251548882317849
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("C4").Value = 0.03740463833779728
$ws.Range("D4").Value = 0.5846394797135872
$ws.Range("E4").Value = 0.07949401871498196
$ws.Range("F4").Value = 7.096297246036784
$ws.Range("G4").Value = 0.002649561795609693
$ws.Range("I4").Value = 6.233440295296163
$ws.Range("L4").Value = 0.2217824092077194
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("C5").Value = 0.03666219987600527
$ws.Range("D5").Value = 0.5786322765804357
$ws.Range("E5").Value = 0.07948633921435722
$ws.Range("F5").Value = 7.010003513989489
$ws.Range("G5").Value = 0.002654464635170271
$ws.Range("I5").Value = 6.164096758339895
$ws.Range("L5").Value = 0.2204598454045623
$ws.Range("N5").Value = 2.293303068605894
$ws.Range("C6").Value = 0.03653879501768387
$ws.Range("D6").Value = 0.5776438682987362
$ws.Range("E6").Value = 0.07948587676149188
$ws.Range("F6").Value = 6.995785942668533
$ws.Range("G6").Value = 0.002655286526118839
$ws.Range("I6").Value = 6.152678665535262
$ws.Range("L6").Value = 0.220243320850301
$ws.Range("N6").Value = 2.2803432614038
$ws.Range("C7").Value = 0.03739463359193707
$ws.Range("D7").Value = 0.5845578522433357
$ws.Range("E7").Value = 0.07949386055604002
$ws.Range("F7").Value = 7.09512594510602
$ws.Range("G7").Value = 0.002649627396373771
$ws.Range("I7").Value = 6.232498609579295
$ws.Range("L7").Value = 0.2217643649056527
$ws.Range("N7").Value = 2.370273851392596
$ws.Range("C8").Value = 0.04115597052486208
$ws.Range("D8").Value = 0.6164245101475387
$ws.Range("E8").Value = 0.0796484072896142
$ws.Range("F8").Value = 7.550163329668123
$ws.Range("G8").Value = 0.002625777846524649
$ws.Range("I8").Value = 6.599122840620311
$ws.Range("L8").Value = 0.228934403684363
$ws.Range("N8").Value = 2.766433886209882
$ws.Range("C9").Value = 0.04850408655065053
$ws.Range("D9").Value = 0.6841645191989869
$ws.Range("E9").Value = 0.08039250942738008
$ws.Range("F9").Value = 8.507258773144315
$ws.Range("G9").Value = 0.002582968601502634
$ws.Range("I9").Value = 7.373821485815938
$ws.Range("L9").Value = 0.2447374443017338
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("C10").Value = 0.05391293722860269
$ws.Range("D10").Value = 0.7375099410671169
$ws.Range("E10").Value = 0.08122049725138325
$ws.Range("F10").Value = 9.254813788461888
$ws.Range("G10").Value = 0.00255385582534831
$ws.Range("I10").Value = 7.981023838710314
$ws.Range("L10").Value = 0.2575078008059961
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("C11").Value = 0.05638207931312422
$ws.Range("D11").Value = 0.7626428374325656
$ws.Range("E11").Value = 0.0816611850509581
$ws.Range("F11").Value = 9.605674321016124
$ws.Range("G11").Value = 0.002541102731681808
$ws.Range("I11").Value = 8.266458780448204
$ws.Range("L11").Value = 0.2635920434573933
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("C12").Value = 0.05731886261308716
$ws.Range("D12").Value = 0.772292150534895
$ws.Range("E12").Value = 0.08183752697458502
$ws.Range("F12").Value = 9.740187184479396
$ws.Range("G12").Value = 0.002536342583612769
$ws.Range("I12").Value = 8.375952612406877
$ws.Range("L12").Value = 0.2659375074848214
$ws.Range("N12").Value = 4.460285735713398
$ws.Range("C13").Value = 0.05711702194739132
$ws.Range("D13").Value = 0.7702080016600235
$ws.Range("E13").Value = 0.08179912348007079
$ws.Range("F13").Value = 9.711142455711865
$ws.Range("G13").Value = 0.002537364712529651
$ws.Range("I13").Value = 8.352307299490803
$ws.Range("L13").Value = 0.2654304920337438
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("C14").Value = 0.05645911036127416
$ws.Range("D14").Value = 0.763434001460837
$ws.Range("E14").Value = 0.08167550163509318
$ws.Range("F14").Value = 9.616707137113394
$ws.Range("G14").Value = 0.002540709732909306
$ws.Range("I14").Value = 8.27543825106244
$ws.Range("L14").Value = 0.2637841629912572
$ws.Range("N14").Value = 4.371107314139238
$ws.Range("C15").Value = 0.05605636766081545
$ws.Range("D15").Value = 0.7593021559002295
$ws.Range("E15").Value = 0.08160101989770396
$ws.Range("F15").Value = 9.559080589523489
$ws.Range("G15").Value = 0.002542767621922312
$ws.Range("I15").Value = 8.228539289177945
$ws.Range("L15").Value = 0.2627812025923362
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("C16").Value = 0.05375179145474362
$ws.Range("D16").Value = 0.7358855031890243
$ws.Range("E16").Value = 0.08119300721998002
$ws.Range("F16").Value = 9.232109634592746
$ws.Range("G16").Value = 0.002554699024554767
$ws.Range("I16").Value = 7.962562223481825
$ws.Range("L16").Value = 0.2571158785252976
$ws.Range("N16").Value = 4.089429168003846
$ws.Range("C17").Value = 0.0523405883303667
$ws.Range("D17").Value = 0.7217471205205186
$ws.Range("E17").Value = 0.08095928151067966
$ws.Range("F17").Value = 9.034355760268909
$ws.Range("G17").Value = 0.002562143223105643
$ws.Range("I17").Value = 7.801810180466703
$ws.Range("L17").Value = 0.2537121625090606
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("C18").Value = 0.05152968032149374
$ws.Range("D18").Value = 0.713696114003227
$ws.Range("E18").Value = 0.08083085946042701
$ws.Range("F18").Value = 8.921622795283042
$ws.Range("G18").Value = 0.002566471176755033
$ws.Range("I18").Value = 7.71021217757459
$ws.Range("L18").Value = 0.2517801903935464
$ws.Range("N18").Value = 3.857331695637072
$ws.Range("C19").Value = 0.05125523879557647
$ws.Range("D19").Value = 0.7109838738061853
$ws.Range("E19").Value = 0.08078840279677024
$ws.Range("F19").Value = 8.883623910112419
$ws.Range("G19").Value = 0.002567944531274395
$ws.Range("I19").Value = 7.679344388810875
$ws.Range("L19").Value = 0.251130423524927
$ws.Range("N19").Value = 3.828614786363971
$ws.Range("C20").Value = 0.05249072895415452
$ws.Range("D20").Value = 0.7232437373290281
$ws.Range("E20").Value = 0.0809835381175894
$ws.Range("F20").Value = 9.055301795545518
$ws.Range("G20").Value = 0.002561345998616787
$ws.Range("I20").Value = 7.818832693362765
$ws.Range("L20").Value = 0.2540718144955747
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("C21").Value = 0.05665230244771635
$ws.Range("D21").Value = 0.7654200440833279
$ws.Range("E21").Value = 0.08171155343018555
$ws.Range("F21").Value = 9.644399505823003
$ws.Range("G21").Value = 0.002539725352916793
$ws.Range("I21").Value = 8.297977735093582
$ws.Range("L21").Value = 0.2642665874595025
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("C22").Value = 0.05938273211287992
$ws.Range("D22").Value = 0.793758085372076
$ws.Range("E22").Value = 0.08224264132876158
$ws.Range("F22").Value = 10.03907326760572
$ws.Range("G22").Value = 0.002525997533223377
$ws.Range("I22").Value = 8.619361651941688
$ws.Range("L22").Value = 0.271172307224262
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("C23").Value = 0.05792429713714853
$ws.Range("D23").Value = 0.7785601720462978
$ws.Range("E23").Value = 0.08195404232722936
$ws.Range("F23").Value = 9.827510510781678
$ws.Range("G23").Value = 0.002533287959621351
$ws.Range("I23").Value = 8.447051655679502
$ws.Range("L23").Value = 0.2674636975786626
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("C24").Value = 0.05242284913732931
$ws.Range("D24").Value = 0.7225668765436524
$ws.Range("E24").Value = 0.08097255319203356
$ws.Range("F24").Value = 9.04582911518645
$ws.Range("G24").Value = 0.002561706273720219
$ws.Range("I24").Value = 7.811134265689134
$ws.Range("L24").Value = 0.2539091385161356
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("C25").Value = 0.04651647950207405
$ws.Range("D25").Value = 0.6652393358801305
$ws.Range("E25").Value = 0.08014275415760252
$ws.Range("F25").Value = 8.240914992244655
$ws.Range("G25").Value = 0.002594133311857528
$ws.Range("I25").Value = 7.157871378860023
$ws.Range("L25").Value = 0.2402659826011018
$ws.Range("N25").Value = 3.331249627311138
